$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    3  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    4  = @(0.7287194209349384, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.034748368925986)
    5  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    6  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    7  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    8  = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795)
    9  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    10 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    11 = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 3.811642989160245)
    12 = @(0.1554434735375247, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569, 0.8605486643198037)
    13 = @(0.3464964993005633, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.896700893398075)
    14 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795)
    15 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    16 = @(1.505614041169197, 86.29678392075563, 3.082599426703578, 6.48142807727062, 97.36642546589903)
    17 = @(0.02258322285507441, 0.004309184025731883, 3.082599426703578, 0.4998867070740569, 3.609378540658442)
    18 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 8.418600821238126)
    19 = @(0.006876353814593728, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1.557080747912106)
    20 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
    21 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    22 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 4.371470058157054)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
